# Auto-generated edit script: refresh market-price derived values in the
# per-class Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Columns H-N are: currentAveragePrice, currentAveragePriceNQ,
# currentAveragePriceHQ, LevePriceNQ, LevePriceHQ, LeveProfitNQ, LeveProfitHQ.

$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(38, 8).Value = 61.125
$ws.Cells.Item(38, 9).Value = 61.125
$ws.Cells.Item(38, 10).Value = 0
$ws.Cells.Item(38, 11).Value = 183.375
$ws.Cells.Item(38, 12).Value = 0
$ws.Cells.Item(38, 13).Value = 188.625
$ws.Cells.Item(100, 8).Value = 45456508
$ws.Cells.Item(100, 9).Value = 90909700
$ws.Cells.Item(100, 10).Value = 3318.7273
$ws.Cells.Item(100, 11).Value = 90909700
$ws.Cells.Item(100, 12).Value = 3318.7273
$ws.Cells.Item(100, 13).Value = -90909159
$ws.Cells.Item(100, 14).Value = -4400.7273
$ws.Cells.Item(129, 8).Value = 295351.6
$ws.Cells.Item(129, 9).Value = 332
$ws.Cells.Item(129, 10).Value = 346217.03
$ws.Cells.Item(129, 11).Value = 996
$ws.Cells.Item(129, 12).Value = 1038651.09
$ws.Cells.Item(129, 13).Value = 4004
$ws.Cells.Item(129, 14).Value = -1048651.09
$ws.Cells.Item(132, 8).Value = 2670.6287
$ws.Cells.Item(132, 9).Value = 3149.7036
$ws.Cells.Item(132, 10).Value = 1053.75
$ws.Cells.Item(132, 11).Value = 9449.110799999999
$ws.Cells.Item(132, 12).Value = 3161.25
$ws.Cells.Item(132, 13).Value = -6919.110799999999
$ws.Cells.Item(132, 14).Value = -8221.25
$ws.Cells.Item(135, 8).Value = 9806716
$ws.Cells.Item(135, 9).Value = 487.39026
$ws.Cells.Item(135, 10).Value = 50012252
$ws.Cells.Item(135, 11).Value = 4386.51234
$ws.Cells.Item(135, 12).Value = 450110268
$ws.Cells.Item(135, 13).Value = -1851.51234
$ws.Cells.Item(135, 14).Value = -450115338
$ws.Cells.Item(137, 8).Value = 1388.92
$ws.Cells.Item(137, 9).Value = 1366.8292
$ws.Cells.Item(137, 10).Value = 1489.5555
$ws.Cells.Item(137, 11).Value = 4100.487599999999
$ws.Cells.Item(137, 12).Value = 4468.666499999999
$ws.Cells.Item(137, 13).Value = -1550.487599999999
$ws.Cells.Item(137, 14).Value = -9568.666499999999
$ws.Cells.Item(138, 8).Value = 19609970
$ws.Cells.Item(138, 9).Value = 28572528
$ws.Cells.Item(138, 10).Value = 4373.5625
$ws.Cells.Item(138, 11).Value = 85717584
$ws.Cells.Item(138, 12).Value = 13120.6875
$ws.Cells.Item(138, 13).Value = -85712444
$ws.Cells.Item(138, 14).Value = -23400.6875
$ws.Cells.Item(141, 8).Value = 1387.2174
$ws.Cells.Item(141, 9).Value = 853
$ws.Cells.Item(141, 10).Value = 3310.4
$ws.Cells.Item(141, 11).Value = 2559
$ws.Cells.Item(141, 12).Value = 9931.200000000001
$ws.Cells.Item(141, 13).Value = 2621
$ws.Cells.Item(141, 14).Value = -20291.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1169.6666
$ws.Cells.Item(2, 9).Value = 1259.625
$ws.Cells.Item(2, 10).Value = 989.75
$ws.Cells.Item(2, 11).Value = 1259.625
$ws.Cells.Item(2, 12).Value = 989.75
$ws.Cells.Item(2, 13).Value = -1146.625
$ws.Cells.Item(2, 14).Value = -1215.75
$ws.Cells.Item(32, 8).Value = 3543.2598
$ws.Cells.Item(32, 9).Value = 3176.75
$ws.Cells.Item(32, 10).Value = 6312.4443
$ws.Cells.Item(32, 11).Value = 3176.75
$ws.Cells.Item(32, 12).Value = 6312.4443
$ws.Cells.Item(32, 13).Value = -2889.75
$ws.Cells.Item(32, 14).Value = -6886.4443
$ws.Cells.Item(61, 8).Value = 347337.06
$ws.Cells.Item(61, 9).Value = 419564.62
$ws.Cells.Item(61, 10).Value = 2249.889
$ws.Cells.Item(61, 11).Value = 419564.62
$ws.Cells.Item(61, 12).Value = 2249.889
$ws.Cells.Item(61, 13).Value = -419352.62
$ws.Cells.Item(61, 14).Value = -2673.889
$ws.Cells.Item(74, 8).Value = 33335336
$ws.Cells.Item(74, 9).Value = 40002150
$ws.Cells.Item(74, 10).Value = 1251.8
$ws.Cells.Item(74, 11).Value = 40002150
$ws.Cells.Item(74, 12).Value = 1251.8
$ws.Cells.Item(74, 13).Value = -40001276
$ws.Cells.Item(74, 14).Value = -2999.8
$ws.Cells.Item(77, 8).Value = 33335336
$ws.Cells.Item(77, 9).Value = 40002150
$ws.Cells.Item(77, 10).Value = 1251.8
$ws.Cells.Item(77, 11).Value = 200010750
$ws.Cells.Item(77, 12).Value = 6259
$ws.Cells.Item(77, 13).Value = -200006382
$ws.Cells.Item(77, 14).Value = -14995
$ws.Cells.Item(97, 8).Value = 2117.7896
$ws.Cells.Item(97, 9).Value = 2033.625
$ws.Cells.Item(97, 10).Value = 2566.6667
$ws.Cells.Item(97, 11).Value = 2033.625
$ws.Cells.Item(97, 12).Value = 2566.6667
$ws.Cells.Item(97, 13).Value = -1537.625
$ws.Cells.Item(102, 8).Value = 1088.4546
$ws.Cells.Item(102, 9).Value = 1088.4546
$ws.Cells.Item(102, 10).Value = 0
$ws.Cells.Item(102, 11).Value = 1088.4546
$ws.Cells.Item(102, 12).Value = 0
$ws.Cells.Item(102, 13).Value = 533.5454
$ws.Cells.Item(116, 8).Value = 1169.6666
$ws.Cells.Item(116, 9).Value = 1259.625
$ws.Cells.Item(116, 10).Value = 989.75
$ws.Cells.Item(116, 11).Value = 1259.625
$ws.Cells.Item(116, 12).Value = 989.75
$ws.Cells.Item(116, 13).Value = 1034.375
$ws.Cells.Item(116, 14).Value = -5577.75
$ws.Cells.Item(132, 8).Value = 19434
$ws.Cells.Item(132, 9).Value = 2185.182
$ws.Cells.Item(132, 10).Value = 73644.57000000001
$ws.Cells.Item(132, 11).Value = 6555.545999999999
$ws.Cells.Item(132, 12).Value = 220933.71
$ws.Cells.Item(132, 13).Value = -4025.545999999999
$ws.Cells.Item(132, 14).Value = -225993.71
$ws.Cells.Item(136, 8).Value = 347337.06
$ws.Cells.Item(136, 9).Value = 419564.62
$ws.Cells.Item(136, 10).Value = 2249.889
$ws.Cells.Item(136, 11).Value = 1258693.86
$ws.Cells.Item(136, 12).Value = 6749.667
$ws.Cells.Item(136, 13).Value = -1256143.86
$ws.Cells.Item(136, 14).Value = -11849.667

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1169.6666
$ws.Cells.Item(3, 9).Value = 1259.625
$ws.Cells.Item(3, 10).Value = 989.75
$ws.Cells.Item(3, 11).Value = 1259.625
$ws.Cells.Item(3, 12).Value = 989.75
$ws.Cells.Item(3, 13).Value = -1145.625
$ws.Cells.Item(3, 14).Value = -1217.75
$ws.Cells.Item(99, 8).Value = 1598.7778
$ws.Cells.Item(99, 9).Value = 1555.5714
$ws.Cells.Item(99, 10).Value = 1750
$ws.Cells.Item(99, 11).Value = 1555.5714
$ws.Cells.Item(99, 12).Value = 1750
$ws.Cells.Item(99, 13).Value = -57.57140000000004
$ws.Cells.Item(99, 14).Value = -4746
$ws.Cells.Item(105, 8).Value = 1847.15
$ws.Cells.Item(105, 9).Value = 1789.5161
$ws.Cells.Item(105, 10).Value = 2045.6666
$ws.Cells.Item(105, 11).Value = 1789.5161
$ws.Cells.Item(105, 12).Value = 2045.6666
$ws.Cells.Item(105, 13).Value = -42.51610000000005
$ws.Cells.Item(105, 14).Value = -5539.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3972.85
$ws.Cells.Item(31, 9).Value = 2460.3845
$ws.Cells.Item(31, 10).Value = 6781.7144
$ws.Cells.Item(31, 11).Value = 2460.3845
$ws.Cells.Item(31, 12).Value = 6781.7144
$ws.Cells.Item(31, 13).Value = -2165.3845
$ws.Cells.Item(31, 14).Value = -7371.7144
$ws.Cells.Item(34, 8).Value = 3972.85
$ws.Cells.Item(34, 9).Value = 2460.3845
$ws.Cells.Item(34, 10).Value = 6781.7144
$ws.Cells.Item(34, 11).Value = 2460.3845
$ws.Cells.Item(34, 12).Value = 6781.7144
$ws.Cells.Item(34, 13).Value = -2258.3845
$ws.Cells.Item(34, 14).Value = -7185.7144
$ws.Cells.Item(58, 8).Value = 13078.523
$ws.Cells.Item(58, 9).Value = 1006.4667
$ws.Cells.Item(58, 10).Value = 43258.668
$ws.Cells.Item(58, 11).Value = 1006.4667
$ws.Cells.Item(58, 12).Value = 43258.668
$ws.Cells.Item(58, 13).Value = -803.4666999999999
$ws.Cells.Item(58, 14).Value = -43664.668
$ws.Cells.Item(86, 8).Value = 7586082.5
$ws.Cells.Item(86, 9).Value = 4671.0557
$ws.Cells.Item(86, 10).Value = 41702436
$ws.Cells.Item(86, 11).Value = 4671.0557
$ws.Cells.Item(86, 12).Value = 41702436
$ws.Cells.Item(86, 13).Value = -3548.0557
$ws.Cells.Item(86, 14).Value = -41704682
$ws.Cells.Item(89, 8).Value = 7586082.5
$ws.Cells.Item(89, 9).Value = 4671.0557
$ws.Cells.Item(89, 10).Value = 41702436
$ws.Cells.Item(89, 11).Value = 23355.2785
$ws.Cells.Item(89, 12).Value = 208512180
$ws.Cells.Item(89, 13).Value = -17739.2785
$ws.Cells.Item(89, 14).Value = -208523412
$ws.Cells.Item(94, 8).Value = 6562.375
$ws.Cells.Item(94, 9).Value = 2500
$ws.Cells.Item(94, 10).Value = 7142.7144
$ws.Cells.Item(94, 11).Value = 2500
$ws.Cells.Item(94, 12).Value = 7142.7144
$ws.Cells.Item(94, 13).Value = -2049
$ws.Cells.Item(94, 14).Value = -8044.7144
$ws.Cells.Item(132, 8).Value = 2478.639
$ws.Cells.Item(132, 9).Value = 1622.0312
$ws.Cells.Item(132, 10).Value = 9331.5
$ws.Cells.Item(132, 11).Value = 4866.0936
$ws.Cells.Item(132, 12).Value = 27994.5
$ws.Cells.Item(132, 13).Value = -2336.0936
$ws.Cells.Item(132, 14).Value = -33054.5
$ws.Cells.Item(134, 8).Value = 740.35895
$ws.Cells.Item(134, 9).Value = 663.1667
$ws.Cells.Item(134, 10).Value = 1666.6666
$ws.Cells.Item(134, 11).Value = 1989.5001
$ws.Cells.Item(134, 12).Value = 4999.9998
$ws.Cells.Item(134, 13).Value = 545.4999
$ws.Cells.Item(136, 8).Value = 13078.523
$ws.Cells.Item(136, 9).Value = 1006.4667
$ws.Cells.Item(136, 10).Value = 43258.668
$ws.Cells.Item(136, 11).Value = 3019.4001
$ws.Cells.Item(136, 12).Value = 129776.004
$ws.Cells.Item(136, 13).Value = -469.4000999999998
$ws.Cells.Item(136, 14).Value = -134876.004

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 6750.067
$ws.Cells.Item(2, 9).Value = 10083.1
$ws.Cells.Item(2, 10).Value = 84
$ws.Cells.Item(2, 11).Value = 60498.60000000001
$ws.Cells.Item(2, 12).Value = 504
$ws.Cells.Item(2, 13).Value = -60385.60000000001
$ws.Cells.Item(2, 14).Value = -730
$ws.Cells.Item(107, 8).Value = 4756.364
$ws.Cells.Item(107, 9).Value = 9270.909
$ws.Cells.Item(107, 10).Value = 241.81818
$ws.Cells.Item(107, 11).Value = 27812.727
$ws.Cells.Item(107, 12).Value = 725.4545400000001
$ws.Cells.Item(107, 13).Value = -25892.727
$ws.Cells.Item(113, 8).Value = 464.5
$ws.Cells.Item(113, 9).Value = 400
$ws.Cells.Item(113, 10).Value = 477.4
$ws.Cells.Item(113, 11).Value = 1200
$ws.Cells.Item(113, 12).Value = 1432.2
$ws.Cells.Item(113, 13).Value = 970
$ws.Cells.Item(113, 14).Value = -5772.2
$ws.Cells.Item(131, 8).Value = 805.1313
$ws.Cells.Item(131, 9).Value = 680
$ws.Cells.Item(131, 10).Value = 814.65216
$ws.Cells.Item(131, 11).Value = 2040
$ws.Cells.Item(131, 12).Value = 2443.95648
$ws.Cells.Item(131, 13).Value = 3000
$ws.Cells.Item(131, 14).Value = -12523.95648
$ws.Cells.Item(132, 8).Value = 780
$ws.Cells.Item(132, 9).Value = 950
$ws.Cells.Item(132, 10).Value = 666.6667
$ws.Cells.Item(132, 11).Value = 8550
$ws.Cells.Item(132, 12).Value = 6000.0003
$ws.Cells.Item(132, 13).Value = -6020
$ws.Cells.Item(132, 14).Value = -11060.0003
$ws.Cells.Item(134, 8).Value = 2712.8635
$ws.Cells.Item(134, 9).Value = 1926
$ws.Cells.Item(134, 10).Value = 4399
$ws.Cells.Item(134, 11).Value = 5778
$ws.Cells.Item(134, 12).Value = 13197
$ws.Cells.Item(134, 13).Value = -708
$ws.Cells.Item(134, 14).Value = -23337

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 4232
$ws.Cells.Item(126, 9).Value = 3223.8262
$ws.Cells.Item(126, 10).Value = 6015.6924
$ws.Cells.Item(126, 11).Value = 9671.4786
$ws.Cells.Item(126, 12).Value = 18047.0772
$ws.Cells.Item(126, 13).Value = -7201.4786
$ws.Cells.Item(126, 14).Value = -22987.0772
$ws.Cells.Item(132, 8).Value = 19035.902
$ws.Cells.Item(132, 9).Value = 2929.6667
$ws.Cells.Item(132, 10).Value = 127753
$ws.Cells.Item(132, 11).Value = 8789.000100000001
$ws.Cells.Item(132, 12).Value = 383259
$ws.Cells.Item(132, 13).Value = -6259.000100000001
$ws.Cells.Item(132, 14).Value = -388319

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 193.25
$ws.Cells.Item(55, 9).Value = 139.25
$ws.Cells.Item(55, 10).Value = 247.25
$ws.Cells.Item(55, 11).Value = 139.25
$ws.Cells.Item(55, 12).Value = 247.25
$ws.Cells.Item(55, 13).Value = 33.75
$ws.Cells.Item(55, 14).Value = -593.25
$ws.Cells.Item(122, 8).Value = 596645
$ws.Cells.Item(122, 9).Value = 1091500.2
$ws.Cells.Item(122, 10).Value = 2818.6
$ws.Cells.Item(122, 11).Value = 3274500.6
$ws.Cells.Item(122, 12).Value = 8455.799999999999
$ws.Cells.Item(122, 13).Value = -3272050.6
$ws.Cells.Item(122, 14).Value = -13355.8
$ws.Cells.Item(136, 8).Value = 918.23334
$ws.Cells.Item(136, 9).Value = 918.23334
$ws.Cells.Item(136, 10).Value = 0
$ws.Cells.Item(136, 11).Value = 2754.70002
$ws.Cells.Item(136, 12).Value = 0
$ws.Cells.Item(136, 13).Value = -204.7000200000002

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 968.96295
$ws.Cells.Item(122, 9).Value = 983.9231
$ws.Cells.Item(122, 10).Value = 580
$ws.Cells.Item(122, 11).Value = 2951.7693
$ws.Cells.Item(122, 12).Value = 1740
$ws.Cells.Item(122, 13).Value = -501.7692999999999
$ws.Cells.Item(122, 14).Value = -6640
$ws.Cells.Item(132, 8).Value = 1057.6086
$ws.Cells.Item(132, 9).Value = 719.97144
$ws.Cells.Item(132, 10).Value = 2131.9092
$ws.Cells.Item(132, 11).Value = 2159.91432
$ws.Cells.Item(132, 12).Value = 6395.7276
$ws.Cells.Item(132, 13).Value = 370.0856800000001
$ws.Cells.Item(132, 14).Value = -11455.7276
$ws.Cells.Item(136, 8).Value = 16951206
$ws.Cells.Item(136, 9).Value = 24391202
$ws.Cells.Item(136, 10).Value = 4550.5557
$ws.Cells.Item(136, 11).Value = 73173606
$ws.Cells.Item(136, 12).Value = 13651.6671
$ws.Cells.Item(136, 13).Value = -73171056
$ws.Cells.Item(136, 14).Value = -18751.6671

# Cells removed entirely in the target revision (column no longer present
# for that row) -- clear their contents so the cell disappears from the XML.

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(38, 14).ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(102, 14).ClearContents()
